$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7, shifting rows 7-9 down to 8-10
$ws.Rows.Item(7).Insert()

# Fill in the new row 7 values
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value = "Ñuble"
$ws.Cells.Item(7, 4).Value = 44466
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = 100112026
$ws.Cells.Item(7, 7).Value = "Haba"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 60
$ws.Cells.Item(7, 11).Value = 11000
$ws.Cells.Item(7, 12).Value = 12000
$ws.Cells.Item(7, 13).Value = 11500
$ws.Cells.Item(7, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(7, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(7, 16).Value = 460
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
